$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 594.75
$ws.Range("I2").Value = 459.66666
$ws.Range("K2").Value = 459.66666
$ws.Range("M2").Value = -346.66666

$ws.Range("H98").Value = 1059.4
$ws.Range("I98").Value = 949.75
$ws.Range("K98").Value = 949.75
$ws.Range("M98").Value = 548.25

$ws.Range("H122").Value = 1059.4
$ws.Range("I122").Value = 949.75
$ws.Range("K122").Value = 2849.25
$ws.Range("M122").Value = -399.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1670.3334
$ws.Range("I110").Value = 1005.5
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1005.5
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 1039.5
$ws.Range("N110").Value = -7090

$ws.Range("H122").Value = 2248
$ws.Range("I122").Value = 1830.6666
$ws.Range("K122").Value = 5491.9998
$ws.Range("M122").Value = -3041.9998

$ws.Range("H132").Value = 1970.5555
$ws.Range("I132").Value = 1851.1765
$ws.Range("K132").Value = 5553.529500000001
$ws.Range("M132").Value = -3023.529500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2490.9092
$ws.Range("I94").Value = 1775
$ws.Range("K94").Value = 1775
$ws.Range("M94").Value = -1324

$ws.Range("H99").Value = 4028.5715
$ws.Range("I99").Value = 4028.5715
$ws.Range("K99").Value = 4028.5715
$ws.Range("M99").Value = -2530.5715

$ws.Range("H105").Value = 1719.8
$ws.Range("I105").Value = 1719.8
$ws.Range("K105").Value = 1719.8
$ws.Range("M105").Value = 27.20000000000005

$ws.Range("H107").Value = 3557
$ws.Range("I107").Value = 3553.4
$ws.Range("J107").Value = 3563
$ws.Range("K107").Value = 3553.4
$ws.Range("L107").Value = 3563
$ws.Range("M107").Value = -1633.4
$ws.Range("N107").Value = -7403

$ws.Range("H134").Value = 3669.6667
$ws.Range("J134").Value = 5014
$ws.Range("L134").Value = 15042
$ws.Range("N134").Value = -20112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1012.6667
$ws.Range("I16").Value = 861.9231
$ws.Range("J16").Value = 1992.5
$ws.Range("K16").Value = 861.9231
$ws.Range("L16").Value = 1992.5
$ws.Range("M16").Value = -574.9231
$ws.Range("N16").Value = -2566.5

$ws.Range("H31").Value = 6250.3335
$ws.Range("I31").Value = 5875.75
$ws.Range("K31").Value = 5875.75
$ws.Range("M31").Value = -5580.75

$ws.Range("H34").Value = 6250.3335
$ws.Range("I34").Value = 5875.75
$ws.Range("K34").Value = 5875.75
$ws.Range("M34").Value = -5673.75

$ws.Range("H107").Value = 1025.6666
$ws.Range("I107").Value = 888.5
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 888.5
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 1031.5
$ws.Range("N107").Value = -5140

$ws.Range("H113").Value = 1012.6667
$ws.Range("I113").Value = 861.9231
$ws.Range("J113").Value = 1992.5
$ws.Range("K113").Value = 861.9231
$ws.Range("L113").Value = 1992.5
$ws.Range("M113").Value = 1308.0769
$ws.Range("N113").Value = -6332.5

$ws.Range("H122").Value = 4734.7144
$ws.Range("I122").Value = 764.6667
$ws.Range("K122").Value = 2294.0001
$ws.Range("M122").Value = 155.9998999999998

$ws.Range("H132").Value = 2028.3
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560

$ws.Range("H134").Value = 965.5
$ws.Range("I134").Value = 917.7143
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 2753.1429
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -218.1428999999998
$ws.Range("N134").Value = -8970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1608.8695
$ws.Range("I4").Value = 1495
$ws.Range("J4").Value = 2018.8
$ws.Range("K4").Value = 4485
$ws.Range("L4").Value = 6056.4
$ws.Range("M4").Value = -4373
$ws.Range("N4").Value = -6280.4

$ws.Range("H7").Value = 167.16667
$ws.Range("I7").Value = 167.16667
$ws.Range("K7").Value = 501.50001
$ws.Range("M7").Value = -389.50001

$ws.Range("H11").Value = 1001
$ws.Range("I11").Value = 1001
$ws.Range("K11").Value = 3003
$ws.Range("M11").Value = -2863

$ws.Range("H132").Value = 986.75
$ws.Range("I132").Value = 1199
$ws.Range("J132").Value = 916
$ws.Range("K132").Value = 10791
$ws.Range("L132").Value = 8244
$ws.Range("M132").Value = -8261
$ws.Range("N132").Value = -13304

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3337
$ws.Range("I80").Value = 3502.5
$ws.Range("J80").Value = 3006
$ws.Range("K80").Value = 3502.5
$ws.Range("L80").Value = 3006
$ws.Range("M80").Value = -2504.5
$ws.Range("N80").Value = -5002

$ws.Range("H83").Value = 3337
$ws.Range("I83").Value = 3502.5
$ws.Range("J83").Value = 3006
$ws.Range("K83").Value = 17512.5
$ws.Range("L83").Value = 15030
$ws.Range("M83").Value = -12520.5
$ws.Range("N83").Value = -25014

$ws.Range("H122").Value = 5096.7617
$ws.Range("I122").Value = 3631.182
$ws.Range("J122").Value = 6708.9
$ws.Range("K122").Value = 10893.546
$ws.Range("L122").Value = 20126.7
$ws.Range("M122").Value = -8443.545999999998
$ws.Range("N122").Value = -25026.7

$ws.Range("H132").Value = 1874.75
$ws.Range("I132").Value = 1874.75
$ws.Range("K132").Value = 5624.25
$ws.Range("M132").Value = -3094.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8874.034
$ws.Range("I7").Value = 8605.444
$ws.Range("K7").Value = 8605.444
$ws.Range("M7").Value = -8493.444

$ws.Range("H82").Value = 3009.2666
$ws.Range("I82").Value = 1770.625
$ws.Range("K82").Value = 1770.625
$ws.Range("M82").Value = -1409.625

$ws.Range("H85").Value = 3009.2666
$ws.Range("I85").Value = 1770.625
$ws.Range("K85").Value = 1770.625
$ws.Range("M85").Value = -522.625

$ws.Range("H126").Value = 8874.034
$ws.Range("I126").Value = 8605.444
$ws.Range("K126").Value = 25816.332
$ws.Range("M126").Value = -23346.332

$ws.Range("H132").Value = 6499.75
$ws.Range("I132").Value = 5666.6665
$ws.Range("K132").Value = 16999.9995
$ws.Range("M132").Value = -14469.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1391.1428
$ws.Range("I122").Value = 935.25
$ws.Range("K122").Value = 2805.75
$ws.Range("M122").Value = -355.75

$ws.Range("H132").Value = 3832.3333
$ws.Range("I132").Value = 3832.3333
$ws.Range("K132").Value = 11496.9999
$ws.Range("M132").Value = -8966.999899999999

$ws.Range("H136").Value = 1241.5333
$ws.Range("I136").Value = 1241.5333
$ws.Range("K136").Value = 3724.5999
$ws.Range("M136").Value = -1174.5999

